$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated S-val data (filtered save games) - values per row, columns B:G
$data = @{
    2 = @(0.06328177979961902, 0.004309184025731883, 0.7127328510149897, 6.48142807727062, 0, 7.261751892110961)
    3 = @(0.3464964993005633, 0.3375848360084654, 0.7127328510149897, 6.48142807727062, 0, 7.878242263594639)
    4 = @(0.3464964993005633, 9.226618575922256, 0.1529057820181812, 6.48142807727062, 0, 16.20744893451162)
    5 = @(0.00006486019690155054, 0.3375848360084654, 0.7127328510149897, 246.9852506941017, 0, 248.035633241322)
    6 = @(0.000009318123435519965, 0.004309184025731883, 0.7127328510149897, 71517.89157740913, 0, 71518.60862876229)
    7 = @(3.182878228561681, 9.226618575922256, 0.7127328510149897, 6.48142807727062, 1, 19.60365773276954)
    8 = @(0.7287194209349384, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 3.594575437922795)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
